# The document ends its Bibliografia section with a paragraph containing
# "1980.". Immediately after it, the site-generated footer block
# (an empty paragraph, a "Ver no Jupiter ..." line and a "(c) 2020 ..."
# copyright line) must be removed, while the paragraph mark that used to
# close the "1980." run/paragraph and everything that follows the footer
# block (the trailing empty paragraph and the page-break paragraph) stay
# untouched.

$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography ("...1980.") so the
# deletion is anchored to content rather than a hard-coded index.
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*1980.*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph ending in '1980.'"
}

# The three paragraphs to remove are the ones right after the anchor:
#   anchorIndex + 1 -> empty paragraph
#   anchorIndex + 2 -> "Ver no Jupiter Salvar em pdf Salvar em docx"
#   anchorIndex + 3 -> "(c) 2020 . Contact: ..." footer line
$startPara = $d.Paragraphs.Item($anchorIndex + 1)
$endPara = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
